$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GUID entry - force Text format first so the leading zeros in
# "000312" are preserved (otherwise COM auto-coerces the numeric-looking
# string into a Number), then strip the format change back off so the
# cell's style stays the default (matches the style of rows 1-3).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "000312"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. baton length of 50 - transforming along z axis now.. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 24-Mar-2023 14:04:17"
